$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.793.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "'1.860.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'244.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "'0.6440"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.50%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.07550"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.83%  "
$ws.Range("D9").Value = "'0.2975"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("D10").Value = "'24.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.20%  "
$ws.Range("D11").Value = "'0.07686"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'1.864.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "'5.052"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "'0.6933"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").Value = "'84.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "'0.000009864"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.83%  "
$ws.Range("D17").Value = "'6.152"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.31%  "
$ws.Range("D18").Value = "'29.804.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("D19").Value = "'2.117.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").Value = "'237.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'7.520"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'158.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("D27").Value = "'8.563"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "'0.06239"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.76%  "
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").Value = "'1.292"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.91%  "
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").Value = "'1.905"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("E35").Value = "  +3.53%  "
$ws.Range("D36").Value = "'0.7300"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "'2.609"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "'2.819"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "'0.01789"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").Value = "'1.216.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").Value = "'6.321"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("D42").Value = "'0.9219"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "'2.028.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("D45").Value = "'102.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "'67.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").Value = "'0.4062"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "'9.187"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'1.678"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.05%  "
$ws.Range("D51").Value = "'0.05785"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
